# Update the "dSF" column (column F) values on the active worksheet.
# This reflects a data repull where the final score-differential figures
# (dSF) were recalculated/refreshed for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 4
    4  = -6
    5  = 2
    6  = -6
    7  = -1
    8  = -3
    9  = 5
    10 = 4
    12 = -5
    13 = 0
    14 = -4
    16 = -2
    17 = -1
    18 = 1
    19 = -3
    20 = -1
    21 = 1
    23 = -3
    24 = 1
    25 = -3
    26 = 1
    27 = 6
    28 = 3
    29 = 3
    30 = 3
    31 = -4
    32 = -3
    33 = -1
    34 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
